$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-18 Wednesday", 2)

$d.Content.Find.Execute("28×25=", $true, $false, $false, $false, $false, $true, 1, $false, "35×16=", 2)
$d.Content.Find.Execute("23×78=", $true, $false, $false, $false, $false, $true, 1, $false, "65×39=", 2)
$d.Content.Find.Execute("50×21=", $true, $false, $false, $false, $false, $true, 1, $false, "97×37=", 2)
$d.Content.Find.Execute("89×53=", $true, $false, $false, $false, $false, $true, 1, $false, "45×84=", 2)
$d.Content.Find.Execute("12×83=", $true, $false, $false, $false, $false, $true, 1, $false, "46×66=", 2)

$d.Content.Find.Execute("13×98=", $true, $false, $false, $false, $false, $true, 1, $false, "98×76=", 2)
$d.Content.Find.Execute("67×70=", $true, $false, $false, $false, $false, $true, 1, $false, "69×18=", 2)
$d.Content.Find.Execute("22×42=", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=", 2)
$d.Content.Find.Execute("83×49=", $true, $false, $false, $false, $false, $true, 1, $false, "38×62=", 2)
$d.Content.Find.Execute("79×85=", $true, $false, $false, $false, $false, $true, 1, $false, "33×71=", 2)

$d.Content.Find.Execute("21×93=", $true, $false, $false, $false, $false, $true, 1, $false, "50×75=", 2)
$d.Content.Find.Execute("58×93=", $true, $false, $false, $false, $false, $true, 1, $false, "61×53=", 2)
$d.Content.Find.Execute("50×13=", $true, $false, $false, $false, $false, $true, 1, $false, "80×12=", 2)
$d.Content.Find.Execute("24×91=", $true, $false, $false, $false, $false, $true, 1, $false, "57×51=", 2)
$d.Content.Find.Execute("98×38=", $true, $false, $false, $false, $false, $true, 1, $false, "16×21=", 2)

$d.Content.Find.Execute("47×26=", $true, $false, $false, $false, $false, $true, 1, $false, "35×63=", 2)
$d.Content.Find.Execute("20×67=", $true, $false, $false, $false, $false, $true, 1, $false, "49×23=", 2)
$d.Content.Find.Execute("62×35=", $true, $false, $false, $false, $false, $true, 1, $false, "28×13=", 2)
$d.Content.Find.Execute("13×68=", $true, $false, $false, $false, $false, $true, 1, $false, "62×47=", 2)
$d.Content.Find.Execute("89×58=", $true, $false, $false, $false, $false, $true, 1, $false, "74×40=", 2)

$d.Content.Find.Execute("75×45=", $true, $false, $false, $false, $false, $true, 1, $false, "63×79=", 2)
$d.Content.Find.Execute("75×14=", $true, $false, $false, $false, $false, $true, 1, $false, "60×92=", 2)
$d.Content.Find.Execute("27×44=", $true, $false, $false, $false, $false, $true, 1, $false, "18×15=", 2)
$d.Content.Find.Execute("63×35=", $true, $false, $false, $false, $false, $true, 1, $false, "52×53=", 2)
$d.Content.Find.Execute("82×15=", $true, $false, $false, $false, $false, $true, 1, $false, "53×38=", 2)
